# Updates crypto price/volume values on Sheet1 to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
$updates = [ordered]@{
    "D2" = "27.115.93"
    "E2" = "  -0.24%  "
    "D3" = "1.897.15"
    "E3" = "  -0.31%  "
    "E4" = "  +0.22%  "
    "D5" = "306.77"
    "E5" = "  +0.10%  "
    "E6" = "  +0.17%  "
    "D7" = "0.5233"
    "E7" = "  -0.47%  "
    "E8" = "  +0.87%  "
    "D9" = "0.07287"
    "E9" = "  +0.38%  "
    "E10" = "  +1.25%  "
    "D11" = "0.9026"
    "E11" = "  +0.48%  "
    "D12" = "0.08175"
    "E12" = "  -2.68%  "
    "D13" = "95.29"
    "E13" = "  +0.70%  "
    "D14" = "1.848.08"
    "E14" = "  -3.04%  "
    "D15" = "5.346"
    "E15" = "  +1.46%  "
    "E16" = "  +0.19%  "
    "E17" = "  +0.52%  "
    "D18" = "14.67"
    "E18" = "  +0.88%  "
    "E19" = "  +0.16%  "
    "D20" = "27.154.00"
    "E20" = "  -0.23%  "
    "D21" = "5.103"
    "E21" = "  +0.88%  "
    "D23" = "6.458"
    "E23" = "  +0.32%  "
    "D24" = "2.349"
    "E24" = "  +3.23%  "
    "D25" = "149.11"
    "E25" = "  +1.49%  "
    "D26" = "18.21"
    "E26" = "  +0.29%  "
    "E27" = "  -0.75%  "
    "D28" = "115.57"
    "D29" = "4.824"
    "E29" = "  +0.64%  "
    "D30" = "4.885"
    "E30" = "  -0.80%  "
    "D31" = "0.09211"
    "E31" = "  -0.81%  "
    "D32" = "0.05039"
    "E32" = "  -0.37%  "
    "D33" = "0.7930"
    "E33" = "  -2.09%  "
    "D34" = "1.221"
    "E34" = "  -0.96%  "
    "D35" = "2.969"
    "E35" = "  +0.60%  "
    "D36" = "3.358"
    "E36" = "  +0.31%  "
    "D37" = "2.644"
    "E37" = "  +1.46%  "
    "D38" = "0.5688"
    "E38" = "  -0.16%  "
    "D39" = "0.01991"
    "E39" = "  -0.01%  "
    "E40" = "  +0.88%  "
    "D41" = "9.039"
    "E41" = "  +0.97%  "
    "D42" = "6.591"
    "E42" = "  -1.03%  "
    "D43" = "116.39"
    "E43" = "  -1.35%  "
    "E44" = "  -0.12%  "
    "D45" = "0.4876"
    "E45" = "  +0.78%  "
    "D46" = "1.003"
    "E46" = "  +0.21%  "
    "D47" = "10.11"
    "E47" = "  -0.35%  "
    "E48" = "  +0.92%  "
    "D49" = "38.38"
    "E49" = "  +2.53%  "
    "D50" = "63.91"
    "E50" = "  +0.43%  "
    "E51" = "  +0.44%  "
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)

    # Price values in column D can look like plain numbers (e.g. "306.77"),
    # which Excel would otherwise silently coerce into a numeric cell and
    # round/retype. Force text entry, then restore the default "Normal"
    # cell style so no stray formatting is left behind (matches source,
    # where these cells carry no explicit style).
    $isPriceColumn = $cellRef.StartsWith("D")
    $looksNumeric = $newValue -match "^-?\d+(\.\d+)?$"

    if ($isPriceColumn -and $looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newValue
    }
}
